$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new blank rows at row 401, shifting existing rows 401:473 down to 410:482
$ws.Range("A401:A409").EntireRow.Insert()

# Populate the newly inserted rows with historical data for 2019-11-18 .. 2019-11-28
$ws.Cells.Item(401, 1).Value = 1574035200
$ws.Cells.Item(401, 2).NumberFormat = "@"
$ws.Cells.Item(401, 2).Value = "2019-11-18"
$ws.Cells.Item(401, 2).Style = "Normal"
$ws.Cells.Item(401, 3).NumberFormat = "@"
$ws.Cells.Item(401, 3).Value = "0198"
$ws.Cells.Item(401, 3).Style = "Normal"
$ws.Cells.Item(401, 4).Value = "GDB"
$ws.Cells.Item(401, 5).Value = 0.44
$ws.Cells.Item(401, 6).Value = 0.475
$ws.Cells.Item(401, 7).Value = 0.44
$ws.Cells.Item(401, 8).Value = 0.47
$ws.Cells.Item(401, 9).Value = 19092500

$ws.Cells.Item(402, 1).Value = 1574121600
$ws.Cells.Item(402, 2).NumberFormat = "@"
$ws.Cells.Item(402, 2).Value = "2019-11-19"
$ws.Cells.Item(402, 2).Style = "Normal"
$ws.Cells.Item(402, 3).NumberFormat = "@"
$ws.Cells.Item(402, 3).Value = "0198"
$ws.Cells.Item(402, 3).Style = "Normal"
$ws.Cells.Item(402, 4).Value = "GDB"
$ws.Cells.Item(402, 5).Value = 0.475
$ws.Cells.Item(402, 6).Value = 0.5
$ws.Cells.Item(402, 7).Value = 0.465
$ws.Cells.Item(402, 8).Value = 0.5
$ws.Cells.Item(402, 9).Value = 10855300

$ws.Cells.Item(403, 1).Value = 1574208000
$ws.Cells.Item(403, 2).NumberFormat = "@"
$ws.Cells.Item(403, 2).Value = "2019-11-20"
$ws.Cells.Item(403, 2).Style = "Normal"
$ws.Cells.Item(403, 3).NumberFormat = "@"
$ws.Cells.Item(403, 3).Value = "0198"
$ws.Cells.Item(403, 3).Style = "Normal"
$ws.Cells.Item(403, 4).Value = "GDB"
$ws.Cells.Item(403, 5).Value = 0.505
$ws.Cells.Item(403, 6).Value = 0.515
$ws.Cells.Item(403, 7).Value = 0.485
$ws.Cells.Item(403, 8).Value = 0.49
$ws.Cells.Item(403, 9).Value = 6985400

$ws.Cells.Item(404, 1).Value = 1574294400
$ws.Cells.Item(404, 2).NumberFormat = "@"
$ws.Cells.Item(404, 2).Value = "2019-11-21"
$ws.Cells.Item(404, 2).Style = "Normal"
$ws.Cells.Item(404, 3).NumberFormat = "@"
$ws.Cells.Item(404, 3).Value = "0198"
$ws.Cells.Item(404, 3).Style = "Normal"
$ws.Cells.Item(404, 4).Value = "GDB"
$ws.Cells.Item(404, 5).Value = 0.495
$ws.Cells.Item(404, 6).Value = 0.495
$ws.Cells.Item(404, 7).Value = 0.46
$ws.Cells.Item(404, 8).Value = 0.46
$ws.Cells.Item(404, 9).Value = 9574100

$ws.Cells.Item(405, 1).Value = 1574380800
$ws.Cells.Item(405, 2).NumberFormat = "@"
$ws.Cells.Item(405, 2).Value = "2019-11-22"
$ws.Cells.Item(405, 2).Style = "Normal"
$ws.Cells.Item(405, 3).NumberFormat = "@"
$ws.Cells.Item(405, 3).Value = "0198"
$ws.Cells.Item(405, 3).Style = "Normal"
$ws.Cells.Item(405, 4).Value = "GDB"
$ws.Cells.Item(405, 5).Value = 0.46
$ws.Cells.Item(405, 6).Value = 0.475
$ws.Cells.Item(405, 7).Value = 0.46
$ws.Cells.Item(405, 8).Value = 0.465
$ws.Cells.Item(405, 9).Value = 2686800

$ws.Cells.Item(406, 1).Value = 1574640000
$ws.Cells.Item(406, 2).NumberFormat = "@"
$ws.Cells.Item(406, 2).Value = "2019-11-25"
$ws.Cells.Item(406, 2).Style = "Normal"
$ws.Cells.Item(406, 3).NumberFormat = "@"
$ws.Cells.Item(406, 3).Value = "0198"
$ws.Cells.Item(406, 3).Style = "Normal"
$ws.Cells.Item(406, 4).Value = "GDB"
$ws.Cells.Item(406, 5).Value = 0.465
$ws.Cells.Item(406, 6).Value = 0.465
$ws.Cells.Item(406, 7).Value = 0.44
$ws.Cells.Item(406, 8).Value = 0.455
$ws.Cells.Item(406, 9).Value = 4517900

$ws.Cells.Item(407, 1).Value = 1574726400
$ws.Cells.Item(407, 2).NumberFormat = "@"
$ws.Cells.Item(407, 2).Value = "2019-11-26"
$ws.Cells.Item(407, 2).Style = "Normal"
$ws.Cells.Item(407, 3).NumberFormat = "@"
$ws.Cells.Item(407, 3).Value = "0198"
$ws.Cells.Item(407, 3).Style = "Normal"
$ws.Cells.Item(407, 4).Value = "GDB"
$ws.Cells.Item(407, 5).Value = 0.455
$ws.Cells.Item(407, 6).Value = 0.48
$ws.Cells.Item(407, 7).Value = 0.45
$ws.Cells.Item(407, 8).Value = 0.46
$ws.Cells.Item(407, 9).Value = 5896400

$ws.Cells.Item(408, 1).Value = 1574812800
$ws.Cells.Item(408, 2).NumberFormat = "@"
$ws.Cells.Item(408, 2).Value = "2019-11-27"
$ws.Cells.Item(408, 2).Style = "Normal"
$ws.Cells.Item(408, 3).NumberFormat = "@"
$ws.Cells.Item(408, 3).Value = "0198"
$ws.Cells.Item(408, 3).Style = "Normal"
$ws.Cells.Item(408, 4).Value = "GDB"
$ws.Cells.Item(408, 5).Value = 0.46
$ws.Cells.Item(408, 6).Value = 0.49
$ws.Cells.Item(408, 7).Value = 0.46
$ws.Cells.Item(408, 8).Value = 0.485
$ws.Cells.Item(408, 9).Value = 6808100

$ws.Cells.Item(409, 1).Value = 1574899200
$ws.Cells.Item(409, 2).NumberFormat = "@"
$ws.Cells.Item(409, 2).Value = "2019-11-28"
$ws.Cells.Item(409, 2).Style = "Normal"
$ws.Cells.Item(409, 3).NumberFormat = "@"
$ws.Cells.Item(409, 3).Value = "0198"
$ws.Cells.Item(409, 3).Style = "Normal"
$ws.Cells.Item(409, 4).Value = "GDB"
$ws.Cells.Item(409, 5).Value = 0.485
$ws.Cells.Item(409, 6).Value = 0.495
$ws.Cells.Item(409, 7).Value = 0.48
$ws.Cells.Item(409, 8).Value = 0.48
$ws.Cells.Item(409, 9).Value = 2838200

